$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 466 (pushes old rows 466..487 down to 468..489)
$ws.Rows.Item(466).Insert()
$ws.Rows.Item(466).Insert()

# New row 466: Coliflor, Primera, week of 2021-11-09 (serial 44509), Region Metropolitana
$ws.Range("A466").Value = 6
$ws.Range("B466").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C466").Value = "Metropolitana"
$ws.Range("D466").Value = 44509
$ws.Range("E466").Value = 13
$ws.Range("F466").Value = 100112008
$ws.Range("G466").Value = "Coliflor"
$ws.Range("H466").Value = "Sin especificar"
$ws.Range("I466").Value = "Primera"
$ws.Range("J466").Value = 15700
$ws.Range("K466").Value = 500
$ws.Range("L466").Value = 600
$ws.Range("M466").Value = 537
$ws.Range("N466").Value = "$/unidad"
$ws.Range("O466").Value = "Región Metropolitana"
$ws.Range("P466").Value = 537
$ws.Range("Q466").Value = 1
$ws.Range("R466").Value = "Hortaliza"

# New row 467: Coliflor, Segunda, same week
$ws.Range("A467").Value = 6
$ws.Range("B467").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C467").Value = "Metropolitana"
$ws.Range("D467").Value = 44509
$ws.Range("E467").Value = 13
$ws.Range("F467").Value = 100112008
$ws.Range("G467").Value = "Coliflor"
$ws.Range("H467").Value = "Sin especificar"
$ws.Range("I467").Value = "Segunda"
$ws.Range("J467").Value = 4900
$ws.Range("K467").Value = 350
$ws.Range("L467").Value = 400
$ws.Range("M467").Value = 383
$ws.Range("N467").Value = "$/unidad"
$ws.Range("O467").Value = "Región Metropolitana"
$ws.Range("P467").Value = 383
$ws.Range("Q467").Value = 1
$ws.Range("R467").Value = "Hortaliza"

# Ensure the date cells keep the date/time number format used throughout column D
$ws.Range("D466:D467").NumberFormat = "YYYY-MM-DD HH:MM:SS"
